$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust column widths to fit the new, wider data
$ws.Columns.Item(1).ColumnWidth = 11.0703125
$ws.Columns.Item(2).ColumnWidth = 15.12890625

# Add the two new data rows
# Force the ID-number columns to be stored as text (shared strings), not numeric values
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A3:A3").NumberFormat = "@"

$ws.Range("A2").Value = "337829999"
$ws.Range("B2").Value = "Tran Van Thanh"
$ws.Range("C2").Value = "nbk-vl"

$ws.Range("A3").Value = "334442222"
$ws.Range("B3").Value = "Nguyen Thi Tho"
$ws.Range("C3").Value = "nbk-qn"
